# Powerpoint writer: consolidate text runs when possible.
# Collapse the "First" / " " / "slide" runs in the title of slide 1,
# and the "Third" / " " / "slide" runs in the title of slide 3,
# into single text runs.

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Text = "x"
$tr1.Text = "First slide"

$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Text = "x"
$tr3.Text = "Third slide"
